$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.676.64'
$ws.Range("E2").Value = '  +3.99%  '
$ws.Range("D3").Value = '3.382.14'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.94'
$ws.Range("E5").Value = '  +6.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.93'
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  +3.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.183'
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.588'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.44'
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000279'
$ws.Range("E12").Value = '  +5.82%  '
$ws.Range("D13").Value = '3.930.23'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '638.49'
$ws.Range("E14").Value = '  +11.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.60'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '68.824.87'
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("D17").Value = '3.393.30'
$ws.Range("E17").Value = '  +2.32%  '
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.02'
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.12'
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.913'
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.96'
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.09'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.72'
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.09'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.85'
$ws.Range("E26").Value = '  +5.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.83'
$ws.Range("E27").Value = '  +4.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.85'
$ws.Range("E28").Value = '  +7.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.70'
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.85'
$ws.Range("E30").Value = '  +1.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '612.09'
$ws.Range("E31").Value = '  +7.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.76'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").Value = '4.007.08'
$ws.Range("E33").Value = '  +7.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.10'
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.45'
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("E38").Value = '  +7.23%  '
$ws.Range("E39").Value = '  +6.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("E40").Value = '  +2.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '33.67'
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("D42").Value = '0.0₃0704'
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.344'
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0424'
$ws.Range("E45").Value = '  +4.00%  '
$ws.Range("E46").Value = '  +2.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.60'
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.38'
$ws.Range("E49").Value = '  +12.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.64'
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.82'
$ws.Range("E51").Value = '  +7.02%  '
